$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from an existing row onto the new row's date cell
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("A7").Value = 42602.582557870373
$ws.Range("B7").Value = "Named"
$ws.Range("C7").Value = 10313
$ws.Range("D7").Value = 6215
$ws.Range("E7").Value = 372
$ws.Range("F7").Value = 61
$ws.Range("G7").Value = 26
$ws.Range("H7").Value = 69
$ws.Range("I7").Value = 29
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 50
